$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing three dated rows down by one (row2->row3, row3->row4, row4->row5)
# so the newest result can take row 2. Row-level formatting (ht/customHeight) for
# rows 3-15 is fixed by position in this sheet, so only cell VALUES move - done
# here by writing literal content rather than relying on a native row-insert
# (which would drag the per-row height along with it).

$ws.Range("A5").Value = "15/6/2025 (Sun)"
$ws.Range("B5").Value = "2 9 3 2`n0 2 8 9`n8 5 2 5`n6 7 4 1"
$ws.Range("C5").Value = "✅ Direct: 12/3547 (0.34%)`n✅ iBet: 12/195 (6.15%)"

$ws.Range("A4").Value = "18/6/2025 (Wed)"
$ws.Range("B4").Value = "2 1 2 1`n3 2 4 7`n0 5 6 9`n5 3 3 8"
$ws.Range("C4").Value = "✅ Direct: 9/3416 (0.26%)`n✅ iBet: 9/188 (4.79%)"

$ws.Range("A3").Value = "21/6/2025 (Sat)"
$ws.Range("B3").Value = "4 6 1 8`n8 1 0 4`n1 5 7 7`n0 2 9 3"
$ws.Range("C3").Value = "✅ Direct: 13/3814 (0.34%)`n✅ iBet: 13/208 (6.25%)"

$ws.Range("A2").Value = "22/6/2025 (Sun)"
$ws.Range("B2").Value = "4 1 3 7`n6 2 5 4`n0 4 2 8`n9 5 6 3"
$ws.Range("C2").Value = "✅ Direct: 11/4144 (0.27%)`n✅ iBet: 11/222 (4.95%)"

# New empty placeholder cell that appears now that row 8 has two dated rows above it
# pushed out of the "3 real entries" block (mirrors the B/C pairing used in rows 5-7)
$ws.Range("C8").Value = 0

# Extend the trailing placeholder block by one row (B9..B33 pattern continues to B34)
$ws.Range("B34").Value = 0
